$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 386-464 (update through 2021-12-08)
$newRows = @(
    @(44460, 0, 3, 42.67425320056899),
    @(44461, 0, 3, 42.67425320056899),
    @(44462, 0, 3, 42.67425320056899),
    @(44463, 1, 2, 28.44950213371266),
    @(44464, 0, 2, 28.44950213371266),
    @(44465, 0, 1, 14.22475106685633),
    @(44466, 1, 2, 28.44950213371266),
    @(44467, 0, 2, 28.44950213371266),
    @(44468, 0, 2, 28.44950213371266),
    @(44469, 0, 2, 28.44950213371266),
    @(44470, 1, 2, 28.44950213371266),
    @(44471, 0, 2, 28.44950213371266),
    @(44472, 0, 2, 28.44950213371266),
    @(44473, 0, 1, 14.22475106685633),
    @(44474, 0, 1, 14.22475106685633),
    @(44475, 0, 1, 14.22475106685633),
    @(44476, 0, 1, 14.22475106685633),
    @(44477, 0, 0, 0),
    @(44478, 0, 0, 0),
    @(44479, 0, 0, 0),
    @(44480, 0, 0, 0),
    @(44481, 0, 0, 0),
    @(44482, 0, 0, 0),
    @(44483, 0, 0, 0),
    @(44484, 0, 0, 0),
    @(44485, 0, 0, 0),
    @(44486, 0, 0, 0),
    @(44487, 0, 0, 0),
    @(44488, 0, 0, 0),
    @(44489, 0, 0, 0),
    @(44490, 0, 0, 0),
    @(44491, 0, 0, 0),
    @(44492, 0, 0, 0),
    @(44493, 0, 0, 0),
    @(44494, 0, 0, 0),
    @(44495, 0, 0, 0),
    @(44496, 0, 0, 0),
    @(44497, 0, 0, 0),
    @(44498, 0, 0, 0),
    @(44499, 0, 0, 0),
    @(44500, 0, 0, 0),
    @(44501, 0, 0, 0),
    @(44502, 1, 1, 14.22475106685633),
    @(44503, 0, 1, 14.22475106685633),
    @(44504, 0, 1, 14.22475106685633),
    @(44505, 0, 1, 14.22475106685633),
    @(44506, 0, 1, 14.22475106685633),
    @(44507, 1, 2, 28.44950213371266),
    @(44508, 1, 3, 42.67425320056899),
    @(44509, 0, 2, 28.44950213371266),
    @(44510, 0, 2, 28.44950213371266),
    @(44511, 0, 2, 28.44950213371266),
    @(44512, 0, 2, 28.44950213371266),
    @(44513, 0, 2, 28.44950213371266),
    @(44514, 1, 2, 28.44950213371266),
    @(44515, 0, 1, 14.22475106685633),
    @(44516, 3, 4, 56.89900426742533),
    @(44517, 0, 4, 56.89900426742533),
    @(44518, 0, 4, 56.89900426742533),
    @(44519, 1, 5, 71.12375533428165),
    @(44520, 0, 5, 71.12375533428165),
    @(44521, 0, 4, 56.89900426742533),
    @(44522, 3, 7, 99.5732574679943),
    @(44523, 1, 5, 71.12375533428165),
    @(44524, 3, 8, 113.7980085348507),
    @(44525, 1, 9, 128.022759601707),
    @(44526, 7, 15, 213.371266002845),
    @(44527, 0, 15, 213.371266002845),
    @(44528, 1, 16, 227.5960170697013),
    @(44529, 1, 14, 199.1465149359886),
    @(44530, 1, 14, 199.1465149359886),
    @(44531, 0, 11, 156.4722617354196),
    @(44532, 2, 12, 170.697012802276),
    @(44533, 5, 10, 142.2475106685633),
    @(44534, 1, 11, 156.4722617354196),
    @(44535, 5, 15, 213.371266002845),
    @(44536, 3, 17, 241.8207681365576),
    @(44537, 6, 22, 312.9445234708393),
    @(44538, 0, 22, 312.9445234708393),
)

# Copy the date-column formatting (style) from the last existing data row (A385)
$ws.Cells.Item(385, 1).Copy() | Out-Null

$startRow = 386
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

